$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the new columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team record (W/L/T) for every data row
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 83   # AD
    $ws.Cells.Item($r, 31).Value = 79   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
